# PROS-10194 - CCRU - Promo tracking KPIs
# Update target/weight figures on the Promo Tracking KPI sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Promo Display" location block
$ws.Range("C2").Value = 60   # Target for PROMO_COMPLIANCE_LOCATION
$ws.Range("D3").Value = 10   # Weight for PROMO_COMPLIANCE_DISPLAY_PRESENCE
$ws.Range("D4").Value = 20   # Weight for PROMO_COMPLIANCE_DISTRIBUTION_TARGET
$ws.Range("D5").Value = 30   # Weight for PROMO_COMPLIANCE_FACINGS_TARGET
$ws.Range("D7").Value = 30   # Weight for PROMO_COMPLIANCE_PRICE_AVAILABILITY_TOTAL
$ws.Range("D8").Value = 0    # Weight for PROMO_COMPLIANCE_PRICE_TARGET

# "Main Shelf" location block
$ws.Range("D11").Value = 20  # Weight for PROMO_COMPLIANCE_DISTRIBUTION_TARGET
$ws.Range("D12").Value = 0   # Weight for PROMO_COMPLIANCE_FACINGS_TARGET
$ws.Range("D15").Value = 30  # Weight for PROMO_COMPLIANCE_PRICE_TARGET

# Leave the cursor where the author last left it before saving.
$ws.Range("B23").Select()
